$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(11, 1).NumberFormat = "@"
$ws.Cells.Item(11, 1).Value = "2025-03-13"
$ws.Cells.Item(11, 1).Style = "Normal"
$ws.Cells.Item(11, 2).Value = 4
$ws.Cells.Item(11, 3).Value = "Amna"
$ws.Cells.Item(11, 4).Value = "01:05:56"
$ws.Cells.Item(11, 5).Value = "01:06:14"

$ws.Cells.Item(12, 1).NumberFormat = "@"
$ws.Cells.Item(12, 1).Value = "2025-03-13"
$ws.Cells.Item(12, 1).Style = "Normal"
$ws.Cells.Item(12, 2).Value = 3
$ws.Cells.Item(12, 3).Value = "nabeel"
$ws.Cells.Item(12, 4).Value = "01:18:49"
$ws.Cells.Item(12, 5).Value = "01:18:55"
